$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. MPCbS sheet: split "natural gas nonpeaker" into "natural gas steam
#    turbine" (existing row) and "natural gas combined cycle" (new row).
# ---------------------------------------------------------------------------
$mpcbs = $wb.Worksheets.Item("MPCbS")

# Rename the existing "natural gas nonpeaker" row to "natural gas steam turbine".
$mpcbs.Range("A3").Value = "natural gas steam turbine"

# Insert a new row right below it for "natural gas combined cycle" - this
# shifts nuclear (and everything after it) down by one row, and Excel will
# automatically re-point the relative formulas (e.g. B11 -> B12) used
# further down the sheet.
$mpcbs.Rows.Item(4).Insert()
$mpcbs.Range("A4").Value = "natural gas combined cycle"
$mpcbs.Range("B4").Formula = "=9*10^12"
$mpcbs.Range("B4").Style = $mpcbs.Range("B3").Style

# Widen column A slightly to fit the new, longer labels.
$mpcbs.Columns.Item(1).ColumnWidth = 26.42578125

# ---------------------------------------------------------------------------
# 2. About sheet: update the "Notes" wording now that coal/gas/nuclear caps
#    are described as "fossil fuels" and there are more than three
#    unconstrained electricity sources.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A27").Value = "Maximums for fossil fuels and nuclear are not imposed, as these power types are unlikely to"
$about.Range("A29").Value = "this limit doesn't come into play for these electricity sources.)"
